$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of profit data appended after running the allocation script on 2025-09-05.
# Format A4 as Text first so the date-like string "09/05/2025" is stored as a
# literal string (matching the existing A2/A3 cells) instead of being
# auto-converted into a date serial number. Resetting the style back to
# "Normal" afterwards keeps the cell's style index the same as its
# unformatted neighbours (no leftover number-format styling).
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "09/05/2025"
$ws.Range("A4").Style = "Normal"

$ws.Range("B4").Value = 0.1196901916565665
$ws.Range("C4").Value = 0.8803098083434335
